# "Generate Report for Archive"
# The localization status report is regenerated: the handoff status text
# moves from "Ready for handoff" to "In Translation" on every sheet that
# tracks it (Overview's per-locale status columns, plus each locale
# sheet's own Status column). The Status columns are then re-sized to fit
# the new, shorter text.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: status columns for each locale (E = zh-cn, F = de-de)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5

# --- zh-cn sheet: Status column (C)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5

# --- de-de sheet: Status column (C)
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
